# Updated cryptos list on Tue Aug 13 17:28:45 UTC 2024 with GitHub Actions
# Refresh price/volume columns (and re-sync a few rows whose coin order
# shifted) for the cryptos table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.930.59"
$ws.Range("E2").Value = "  +2.56%  "

$ws.Range("D3").Value = "2.690.59"
$ws.Range("E3").Value = "  +2.00%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.12"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.68"
$ws.Range("E6").Value = "  +1.27%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  +1.06%  "

$ws.Range("D9").Value = "2.710.88"
$ws.Range("E9").Value = "  +1.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("E12").Value = "  +0.68%  "

$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").Value = "3.163.71"
$ws.Range("E14").Value = "  +1.92%  "

$ws.Range("D15").Value = "60.967.27"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").Value = "2.862.25"
$ws.Range("E16").Value = "  +7.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.47"
$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.27"
$ws.Range("E19").Value = "  +2.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.58"
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("E22").Value = "  +2.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.68"
$ws.Range("E24").Value = "  +1.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.424"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.169"
$ws.Range("E26").Value = "  +4.35%  "

$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("D28").Value = "0.0₃0828"
$ws.Range("E28").Value = "  +1.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.29"
$ws.Range("E29").Value = "  +1.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.83"
$ws.Range("E30").Value = "  +5.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.60"
$ws.Range("E32").Value = "  +0.78%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.15"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.37"
$ws.Range("E34").Value = "  +0.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.21"
$ws.Range("E35").Value = "  +3.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.955"
$ws.Range("E36").Value = "  -8.21%  "

$ws.Range("E37").Value = "  +4.23%  "

$ws.Range("E38").Value = "  +10.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.884"
$ws.Range("E39").Value = "  +1.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.88"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.54"
$ws.Range("E42").Value = "  -0.73%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.615"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0993"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.149.87"
$ws.Range("E45").Value = "  +8.00%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.00"
$ws.Range("E46").Value = "  +1.17%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.94"
$ws.Range("E48").Value = "  +4.03%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0544"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.33"
$ws.Range("E51").Value = "  +3.78%  "
